$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 941.8333
$ws.Range("J58").Value = 1599.5
$ws.Range("L58").Value = 4798.5
$ws.Range("N58").Value = -5098.5

$ws.Range("H127").Value = 1850.1305
$ws.Range("I127").Value = 1773.6154
$ws.Range("K127").Value = 5320.8462
$ws.Range("M127").Value = -360.8462

$ws.Range("H132").Value = 7247682.5
$ws.Range("I132").Value = 8334545
$ws.Range("K132").Value = 25003635
$ws.Range("M132").Value = -25001105

$ws.Range("H137").Value = 1422.9333
$ws.Range("I137").Value = 1016.8182
$ws.Range("J137").Value = 2539.75
$ws.Range("K137").Value = 3050.4546
$ws.Range("L137").Value = 7619.25
$ws.Range("M137").Value = -500.4546
$ws.Range("N137").Value = -12719.25

$ws.Range("H138").Value = 1499.1356
$ws.Range("I138").Value = 1488.1632
$ws.Range("J138").Value = 1552.9
$ws.Range("K138").Value = 4464.4896
$ws.Range("L138").Value = 4658.700000000001
$ws.Range("M138").Value = 675.5104000000001
$ws.Range("N138").Value = -14938.7

$ws.Range("H139").Value = 48868.145
$ws.Range("J139").Value = 48868.145
$ws.Range("L139").Value = 48868.145
$ws.Range("N139").Value = -59148.145

$ws.Range("H141").Value = 637748.1
$ws.Range("I141").Value = 757468.4
$ws.Range("J141").Value = 4941.143
$ws.Range("K141").Value = 2272405.2
$ws.Range("L141").Value = 14823.429
$ws.Range("M141").Value = -2267225.2
$ws.Range("N141").Value = -25183.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1240.3489
$ws.Range("I74").Value = 921.4516
$ws.Range("J74").Value = 2064.1667
$ws.Range("K74").Value = 921.4516
$ws.Range("L74").Value = 2064.1667
$ws.Range("M74").Value = -47.45159999999998
$ws.Range("N74").Value = -3812.1667

$ws.Range("H77").Value = 1240.3489
$ws.Range("I77").Value = 921.4516
$ws.Range("J77").Value = 2064.1667
$ws.Range("K77").Value = 4607.258
$ws.Range("L77").Value = 10320.8335
$ws.Range("M77").Value = -239.2579999999998
$ws.Range("N77").Value = -19056.8335

$ws.Range("H132").Value = 1362.3334
$ws.Range("I132").Value = 1010.7368
$ws.Range("J132").Value = 2197.375
$ws.Range("K132").Value = 3032.2104
$ws.Range("L132").Value = 6592.125
$ws.Range("M132").Value = -502.2103999999999
$ws.Range("N132").Value = -11652.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2647345.2
$ws.Range("I31").Value = 6494902
$ws.Range("K31").Value = 6494902
$ws.Range("M31").Value = -6494607

$ws.Range("H34").Value = 2647345.2
$ws.Range("I34").Value = 6494902
$ws.Range("K34").Value = 6494902
$ws.Range("M34").Value = -6494700

$ws.Range("H41").Value = 26162.666
$ws.Range("J41").Value = 26162.666
$ws.Range("L41").Value = 26162.666
$ws.Range("N41").Value = -27018.666

$ws.Range("H50").Value = 8500
$ws.Range("I50").Value = 2000
$ws.Range("K50").Value = 2000
$ws.Range("M50").Value = -1375

$ws.Range("H51").Value = 31000
$ws.Range("J51").Value = 31000
$ws.Range("L51").Value = 31000
$ws.Range("N51").Value = -32472

$ws.Range("H60").Value = 20324.75
$ws.Range("J60").Value = 20324.75
$ws.Range("L60").Value = 20324.75
$ws.Range("N60").Value = -21346.75

$ws.Range("H61").Value = 31000
$ws.Range("J61").Value = 31000
$ws.Range("L61").Value = 31000
$ws.Range("N61").Value = -31696

$ws.Range("H132").Value = 2327.8262
$ws.Range("I132").Value = 1225
$ws.Range("K132").Value = 3675
$ws.Range("M132").Value = -1145

$ws.Range("H134").Value = 1368.491
$ws.Range("I134").Value = 1212.4
$ws.Range("K134").Value = 3637.2
$ws.Range("M134").Value = -1102.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7622.5713
$ws.Range("I56").Value = 7622.5713
$ws.Range("K56").Value = 7622.5713
$ws.Range("M56").Value = -7092.5713

$ws.Range("H124").Value = 3030
$ws.Range("I124").Value = 3030
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 9090
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -4180
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 1780
$ws.Range("I125").Value = 1780
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 5340
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -420
$ws.Range("N125").ClearContents()

$ws.Range("H129").Value = 56481.848
$ws.Range("I129").Value = 685.2857
$ws.Range("J129").Value = 121577.836
$ws.Range("K129").Value = 2055.8571
$ws.Range("L129").Value = 364733.508
$ws.Range("M129").Value = 2944.1429
$ws.Range("N129").Value = -374733.508

$ws.Range("H131").Value = 814.0617999999999
$ws.Range("J131").Value = 887.8929000000001
$ws.Range("L131").Value = 2663.6787
$ws.Range("N131").Value = -12743.6787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14973.714
$ws.Range("I70").Value = 22979
$ws.Range("K70").Value = 22979
$ws.Range("M70").Value = -22709

$ws.Range("H73").Value = 14973.714
$ws.Range("I73").Value = 22979
$ws.Range("K73").Value = 22979
$ws.Range("M73").Value = -22043

$ws.Range("H132").Value = 1168202.6
$ws.Range("I132").Value = 1924932.1
$ws.Range("K132").Value = 5774796.300000001
$ws.Range("M132").Value = -5772266.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10002.471
$ws.Range("I40").Value = 10095.667
$ws.Range("J40").Value = 9778.799999999999
$ws.Range("K40").Value = 10095.667
$ws.Range("L40").Value = 9778.799999999999
$ws.Range("M40").Value = -9959.666999999999
$ws.Range("N40").Value = -10050.8

$ws.Range("H122").Value = 5825.5454
$ws.Range("I122").Value = 8569.25
$ws.Range("J122").Value = 4257.7144
$ws.Range("K122").Value = 25707.75
$ws.Range("L122").Value = 12773.1432
$ws.Range("M122").Value = -23257.75
$ws.Range("N122").Value = -17673.1432

$ws.Range("H132").Value = 1458.1786
$ws.Range("I132").Value = 984.5306399999999
$ws.Range("K132").Value = 2953.59192
$ws.Range("M132").Value = -423.5919199999998

$ws.Range("H136").Value = 1675.1578
$ws.Range("J136").Value = 4907.7
$ws.Range("L136").Value = 14723.1
$ws.Range("N136").Value = -19823.1

$ws.Range("H141").Value = 49350
$ws.Range("J141").Value = 49350
$ws.Range("L141").Value = 49350
$ws.Range("N141").Value = -59710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 12665.909
$ws.Range("I18").Value = 12666.5
$ws.Range("J18").Value = 12665.777
$ws.Range("K18").Value = 12666.5
$ws.Range("L18").Value = 12665.777
$ws.Range("M18").Value = -12493.5
$ws.Range("N18").Value = -13011.777

$ws.Range("H70").Value = 29776.25
$ws.Range("J70").Value = 29776.25
$ws.Range("L70").Value = 29776.25
$ws.Range("N70").Value = -30406.25

$ws.Range("H73").Value = 29776.25
$ws.Range("J73").Value = 29776.25
$ws.Range("L73").Value = 29776.25
$ws.Range("N73").Value = -31960.25

$ws.Range("H123").Value = 41866.668
$ws.Range("J123").Value = 41866.668
$ws.Range("L123").Value = 41866.668
$ws.Range("N123").Value = -51666.668

$ws.Range("H132").Value = 1281.137
$ws.Range("I132").Value = 838.0847
$ws.Range("J132").Value = 3148.2856
$ws.Range("K132").Value = 2514.2541
$ws.Range("L132").Value = 9444.856800000001
$ws.Range("M132").Value = 15.74589999999989
$ws.Range("N132").Value = -14504.8568
